$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '40.721.86'
$ws.Range('D3').Value = '2.185.65'
$ws.Range('E3').Value = '  -7.52%  '
$ws.Range('E4').Value = '  +0.04%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '240.99'
$ws.Range('E5').Value = '  +0.25%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '0.620'
$ws.Range('E6').Value = '  -7.87%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '68.35'
$ws.Range('E7').Value = '  -7.85%  '
$ws.Range('E8').Value = '  +0.13%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.536'
$ws.Range('E9').Value = '  -13.05%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.0939'
$ws.Range('E10').Value = '  -8.59%  '
$ws.Range('B11').Value = 'OKB'
$ws.Range('C11').Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '57.48'
$ws.Range('E11').Value = '  -5.27%  '
$ws.Range('B12').Value = 'Avalanche'
$ws.Range('C12').Value = 'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax'
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '36.08'
$ws.Range('E12').Value = '  -4.52%  '
$ws.Range('E13').Value = '  -4.59%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '6.55'
$ws.Range('E14').Value = '  -10.48%  '
$ws.Range('D15').Value = '2.511.93'
$ws.Range('E15').Value = '  -7.45%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '14.58'
$ws.Range('E16').Value = '  -11.01%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '0.828'
$ws.Range('E17').Value = '  -10.27%  '
$ws.Range('D18').Value = '2.187.89'
$ws.Range('E18').Value = '  -7.67%  '
$ws.Range('D19').Value = '40.638.34'
$ws.Range('E19').Value = '  -7.33%  '
$ws.Range('D20').Value = '0.0₃0936'
$ws.Range('E20').Value = '  -9.65%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '72.23'
$ws.Range('E21').Value = '  -7.45%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '6.02'
$ws.Range('E22').Value = '  -8.45%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '228.28'
$ws.Range('E23').Value = '  -9.30%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '2.01'
$ws.Range('E25').Value = '  +0.15%  '
$ws.Range('E26').Value = '  -5.21%  '
$ws.Range('E27').Value = '  -4.76%  '
$ws.Range('E29').Value = '  -8.80%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '168.33'
$ws.Range('E30').Value = '  -4.25%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '20.11'
$ws.Range('E31').Value = '  -10.11%  '
$ws.Range('E32').Value = '  -10.08%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '0.122'
$ws.Range('E33').Value = '  -8.73%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '0.0696'
$ws.Range('E34').Value = '  -7.76%  '
$ws.Range('E35').Value = '  -5.77%  '
$ws.Range('E36').Value = '  -10.78%  '
$ws.Range('E37').Value = '  -0.34%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '23.44'
$ws.Range('E38').Value = '  +14.76%  '
$ws.Range('E39').Value = '  -6.92%  '
$ws.Range('E40').Value = '  -4.32%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '5.74'
$ws.Range('E41').Value = '  -13.64%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '61.85'
$ws.Range('E42').Value = '  -4.95%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '4.82'
$ws.Range('E43').Value = '  -11.40%  '
$ws.Range('E44').Value = '  -6.03%  '
$ws.Range('E45').Value = '  -7.73%  '
$ws.Range('E46').Value = '  +0.01%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '0.0974'
$ws.Range('E47').Value = '  -9.05%  '
$ws.Range('E48').Value = '  +0.05%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '10.22'
$ws.Range('E49').Value = '  +5.32%  '
$ws.Range('E50').Value = '  -7.00%  '
$ws.Range('E51').Value = '  -6.72%  '
